# Sprint3 task table: add a "gastou hoje" (time spent today) tracking
# column (F) and record who worked on the first few tasks and how long.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new values in the same order they were originally entered so the
# shared-string table comes out in a matching order.
$ws.Range("E3").Value = "Jader"
$ws.Range("E4").Value = "Ana e Beatriz"
$ws.Range("E2").Value = "Marco e Vinícius"
$ws.Range("F1").Value = "gastou hoje"
$ws.Range("F3").Value = "50min"
$ws.Range("F4").Value = "50min"
$ws.Range("F2").Value = "1h"

# Match the existing header formatting (centered + accent fill) for the new
# header cell, and the existing centered data-row formatting for the new
# data cells, by copying formats from neighboring cells already styled that
# way.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("E2:F4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New column width to match the exported workbook.
$ws.Columns.Item(6).ColumnWidth = 12.109375

# View tweaks captured in the diff: zoom in a bit and leave the selection on
# the next empty time-tracking cell.
$excel.ActiveWindow.Zoom = 135
$ws.Range("G10").Select()
